$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values, repulled/recalculated from source data.
$ws.Range("F12").Value = -8
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = -6
$ws.Range("F19").Value = -5
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = -2
$ws.Range("F28").Value = -7
$ws.Range("F29").Value = -8
$ws.Range("F30").Value = -4
$ws.Range("F31").Value = -3
$ws.Range("F36").Value = 2
$ws.Range("F38").Value = 0
